$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1. Merge the two runs of the "year selector" bullet into one run ---
# Current paragraph (4th bullet) is built from two <a:r> runs:
#   "년도를 선택하면 "                (no dirty attr)
#   "해당 년도의 휴일을 모두 보여주기"   (dirty="0")
# Target: a single run containing both pieces of text, keeping the
# formatting (dirty="0") of the second run.
$shpText = $s.Shapes.Item(3)
$tf = $shpText.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(4, 1)

$firstRunLen = 9   # length of "년도를 선택하면 "
$secondRunLen = 18 # length of "해당 년도의 휴일을 모두 보여주기"

# Rewrite the second run so it carries the full combined text (this keeps
# that run's rPr, i.e. dirty="0"), then clear out the now-redundant first
# run so only a single run remains.
$secondRun = $para.Characters($firstRunLen + 1, $secondRunLen)
$secondRun.Text = "년도를 선택하면 해당 년도의 휴일을 모두 보여주기"

$firstRun = $para.Characters(1, $firstRunLen)
$firstRun.Text = ""

# The shape auto-fits its height to the text (a:spAutoFit); editing the
# run text makes the host re-flow/re-measure it, which would otherwise
# leave a stray height change that isn't part of this edit. Restore the
# original height explicitly (same pt value used by spAutoFit == exact
# original EMU once rounded through the single-precision Height setter).
$shpText.Height = 396.4845275878906

# --- 2. Nudge the red rectangle (id=3) to its new position ---
# Target off = (1652954, 1103140) EMU; size (ext) stays the same.
# Shape.Left/Top are single-precision (like real PowerPoint), so the
# literals below are chosen to round-trip to the exact EMU values.
$shpRect = $s.Shapes.Item(4)
$shpRect.Left = 130.15386962890625
$shpRect.Top = 86.86141967773438

# --- 3. Add a fade-in "appear" entrance animation for that rectangle ---
$seq = $s.TimeLine.MainSequence
$seq.AddEffect($shpRect, 10, 0, 1) | Out-Null
